$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 214.42857
$ws.Range("I2").Value = 197.45454
$ws.Range("J2").Value = 276.66666
$ws.Range("K2").Value = 197.45454
$ws.Range("L2").Value = 276.66666
$ws.Range("M2").Value = -84.45454000000001
$ws.Range("N2").Value = -502.66666

$ws.Range("H9").Value = 126.454544
$ws.Range("I9").Value = 63.142857
$ws.Range("J9").Value = 237.25
$ws.Range("K9").Value = 63.142857
$ws.Range("L9").Value = 237.25
$ws.Range("M9").Value = 105.857143
$ws.Range("N9").Value = -575.25

$ws.Range("H19").Value = 1016.96155
$ws.Range("I19").Value = 718.4737
$ws.Range("J19").Value = 1827.1428
$ws.Range("K19").Value = 718.4737
$ws.Range("L19").Value = 1827.1428
$ws.Range("M19").Value = -543.4737
$ws.Range("N19").Value = -2177.1428

$ws.Range("H31").Value = 725.75
$ws.Range("I31").Value = 725.75
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2177.25
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1947.25
$ws.Range("N31").Value = $null

$ws.Range("H58").Value = 1159.625
$ws.Range("J58").Value = 2899.6667
$ws.Range("L58").Value = 8699.000100000001
$ws.Range("N58").Value = -8999.000100000001

$ws.Range("H70").Value = 2362.2104
$ws.Range("J70").Value = 7725
$ws.Range("L70").Value = 23175
$ws.Range("N70").Value = -23715

$ws.Range("H73").Value = 2362.2104
$ws.Range("J73").Value = 7725
$ws.Range("L73").Value = 23175
$ws.Range("N73").Value = -25047

$ws.Range("H113").Value = 8776.360000000001
$ws.Range("I113").Value = 3676.125
$ws.Range("K113").Value = 3676.125
$ws.Range("M113").Value = -422.125

$ws.Range("H129").Value = 909.1316
$ws.Range("J129").Value = 937.7714
$ws.Range("L129").Value = 2813.3142
$ws.Range("N129").Value = -12813.3142

$ws.Range("H132").Value = 13971.753
$ws.Range("I132").Value = 14811.098
$ws.Range("J132").Value = 1885.2
$ws.Range("K132").Value = 44433.294
$ws.Range("L132").Value = 5655.6
$ws.Range("M132").Value = -41903.294
$ws.Range("N132").Value = -10715.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14753.04
$ws.Range("I32").Value = 12492.091
$ws.Range("K32").Value = 12492.091
$ws.Range("M32").Value = -12205.091

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 354.44446
$ws.Range("I22").Value = 241.42857
$ws.Range("J22").Value = 750
$ws.Range("K22").Value = 241.42857
$ws.Range("L22").Value = 750
$ws.Range("M22").Value = -68.42857000000001
$ws.Range("N22").Value = -1096

$ws.Range("H52").Value = 22711.6
$ws.Range("J52").Value = 22711.6
$ws.Range("L52").Value = 22711.6
$ws.Range("N52").Value = -23237.6

$ws.Range("H105").Value = 2577.75
$ws.Range("I105").Value = 2422
$ws.Range("J105").Value = 3045
$ws.Range("K105").Value = 2422
$ws.Range("L105").Value = 3045
$ws.Range("M105").Value = -675
$ws.Range("N105").Value = -6539

$ws.Range("H107").Value = 425.125
$ws.Range("I107").Value = 400.16666
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 400.16666
$ws.Range("L107").Value = 500
$ws.Range("M107").Value = 1519.83334
$ws.Range("N107").Value = -4340

$ws.Range("H121").Value = 22711.6
$ws.Range("J121").Value = 22711.6
$ws.Range("L121").Value = 22711.6
$ws.Range("N121").Value = -26205.6

$ws.Range("H132").Value = 46457.5
$ws.Range("J132").Value = 46457.5
$ws.Range("L132").Value = 46457.5
$ws.Range("N132").Value = -56577.5

$ws.Range("H134").Value = 2152.6858
$ws.Range("I134").Value = 1679.4546
$ws.Range("J134").Value = 2953.5386
$ws.Range("K134").Value = 5038.3638
$ws.Range("L134").Value = 8860.6158
$ws.Range("M134").Value = -2503.3638
$ws.Range("N134").Value = -13930.6158

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1988.091
$ws.Range("I105").Value = 2552.8572
$ws.Range("J105").Value = 999.75
$ws.Range("K105").Value = 2552.8572
$ws.Range("L105").Value = 999.75
$ws.Range("M105").Value = -805.8571999999999
$ws.Range("N105").Value = -4493.75

$ws.Range("H122").Value = 1224284.9
$ws.Range("I122").Value = 200970
$ws.Range("J122").Value = 2503428.5
$ws.Range("K122").Value = 602910
$ws.Range("L122").Value = 7510285.5
$ws.Range("M122").Value = -600460
$ws.Range("N122").Value = -7515185.5

$ws.Range("H134").Value = 1841.8445
$ws.Range("I134").Value = 1699.225
$ws.Range("J134").Value = 2982.8
$ws.Range("K134").Value = 5097.674999999999
$ws.Range("L134").Value = 8948.400000000001
$ws.Range("M134").Value = -2562.674999999999
$ws.Range("N134").Value = -14018.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 499.8
$ws.Range("I46").Value = 125
$ws.Range("J46").Value = 1999
$ws.Range("K46").Value = 375
$ws.Range("L46").Value = 5997
$ws.Range("M46").Value = -284
$ws.Range("N46").Value = -6179

$ws.Range("H107").Value = 401128.1
$ws.Range("I107").Value = 991.82355
$ws.Range("J107").Value = 801264.4
$ws.Range("K107").Value = 2975.47065
$ws.Range("L107").Value = 2403793.2
$ws.Range("M107").Value = -1055.47065
$ws.Range("N107").Value = -2407633.2

$ws.Range("H122").Value = 960.0278
$ws.Range("I122").Value = 594.5
$ws.Range("J122").Value = 1910.4
$ws.Range("K122").Value = 5350.5
$ws.Range("L122").Value = 17193.6
$ws.Range("M122").Value = -2900.5
$ws.Range("N122").Value = -22093.6

$ws.Range("H131").Value = 870.34
$ws.Range("I131").Value = 489.23077
$ws.Range("J131").Value = 927.2873499999999
$ws.Range("K131").Value = 1467.69231
$ws.Range("L131").Value = 2781.86205
$ws.Range("M131").Value = 3572.30769
$ws.Range("N131").Value = -12861.86205

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2252.182
$ws.Range("I113").Value = 1832
$ws.Range("J113").Value = 2602.3333
$ws.Range("K113").Value = 1832
$ws.Range("L113").Value = 2602.3333
$ws.Range("M113").Value = 338
$ws.Range("N113").Value = -6942.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").Value = $null

$ws.Range("H46").Value = 1801.7858
$ws.Range("I46").Value = 1713.8889
$ws.Range("J46").Value = 1960
$ws.Range("K46").Value = 1713.8889
$ws.Range("L46").Value = 1960
$ws.Range("M46").Value = -1525.8889
$ws.Range("N46").Value = -2336

$ws.Range("H55").Value = 375.33334
$ws.Range("I55").Value = 250
$ws.Range("J55").Value = 438
$ws.Range("K55").Value = 250
$ws.Range("L55").Value = 438
$ws.Range("M55").Value = -77
$ws.Range("N55").Value = -784

$ws.Range("H61").Value = 1860
$ws.Range("I61").Value = 1860
$ws.Range("K61").Value = 1860
$ws.Range("M61").Value = -1658

$ws.Range("H113").Value = 1860
$ws.Range("I113").Value = 1860
$ws.Range("K113").Value = 1860
$ws.Range("M113").Value = 310

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 697.6923
$ws.Range("I107").Value = 706.7
$ws.Range("J107").Value = 667.6667
$ws.Range("K107").Value = 2120.1
$ws.Range("L107").Value = 2003.0001
$ws.Range("M107").Value = -200.1000000000004
$ws.Range("N107").Value = -5843.0001
